$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.599999999999995
$ws.Range("C7").Value = -13.0038
$ws.Range("B9").Value = 6.634899999999993
$ws.Range("C12").Value = -10.8122
$ws.Range("C14").Value = -13.1999
$ws.Range("D15").Value = -8.723100000000001
$ws.Range("B18").Value = 6.096099999999999
$ws.Range("B20").Value = 8.976800000000003
$ws.Range("C26").Value = -12.38789999999999
$ws.Range("B27").Value = 5.425599999999999
$ws.Range("C27").Value = -12.6933
$ws.Range("C29").Value = -11.13780000000001
$ws.Range("D33").Value = -7.722999999999998
$ws.Range("B35").Value = 8.336700000000006
$ws.Range("D35").Value = -8.422299999999996
$ws.Range("C37").Value = -13.7296
$ws.Range("C38").Value = -13.3998
$ws.Range("D38").Value = -8.218699999999998
$ws.Range("D43").Value = -8.088100000000004
$ws.Range("D44").Value = -7.504300000000002
$ws.Range("D47").Value = -7.599900000000003
$ws.Range("C51").Value = -12.2693
$ws.Range("D51").Value = -7.942300000000001
$ws.Range("C52").Value = -11.3925
$ws.Range("C55").Value = -14.0238
$ws.Range("D57").Value = -8.161199999999999
$ws.Range("D63").Value = -7.808799999999998
$ws.Range("B69").Value = 6.458999999999995
$ws.Range("C69").Value = -11.4596
$ws.Range("C70").Value = -12.47350000000001
$ws.Range("D70").Value = -8.299999999999999
$ws.Range("B76").Value = 4.7408
$ws.Range("B78").Value = 9.752899999999999
$ws.Range("C81").Value = -12.79060000000001
$ws.Range("B82").Value = 5.783000000000002
$ws.Range("B83").Value = 6.398399999999999
$ws.Range("C83").Value = -13.6103
$ws.Range("D88").Value = -7.811399999999995
$ws.Range("B93").Value = 6.5463
$ws.Range("D99").Value = -7.567199999999997
$ws.Range("C102").Value = -13.4714
